$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Requisitos" list (rows 23-25, columns B and C) needs to be reordered:
# the "Indicação de Conjunto" entry (LOM3229) moves from the first position
# to the last position, while the two "Requisito" entries shift up.

$line1 = "LOB1021 -  Física IV  (Requisito)`n"
$line2 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$line3 = "LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)`n"

$ws.Range("B23").Value = $line1
$ws.Range("C23").Value = $line1

$ws.Range("B24").Value = $line2
$ws.Range("C24").Value = $line2

$ws.Range("B25").Value = $line3
$ws.Range("C25").Value = $line3
